$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (merged header numbering) value updates
$ws.Range("B2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("H2").Value = 2

# Row 4 (data row) value updates
$ws.Range("B4").Value = 0.7041217312852566
$ws.Range("C4").Value = -1
$ws.Range("D4").Value = 0.4082434625705131
$ws.Range("E4").Value = 0.5497199333020351
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.09943986660407012
$ws.Range("H4").Value = 0.7741065466434253
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0.5482130932868505
$ws.Range("K4").Value = 0
